# Reorder rows 2-7 of the historical-distance sheet (new JSON source was
# merged in a different sort order) while keeping each record's
# title / timestamp / historical distance / time bucket / uri (+ hyperlink)
# together as one unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the six existing records (rows 2..7) before mutating anything.
$records = @()
for ($r = 2; $r -le 7; $r++) {
    $records += [PSCustomObject]@{
        Title     = $ws.Cells.Item($r, 1).Value2
        Timestamp = $ws.Cells.Item($r, 2).Value2
        Distance  = $ws.Cells.Item($r, 3).Value2
        Bucket    = $ws.Cells.Item($r, 4).Value2
        Uri       = $ws.Cells.Item($r, 5).Value2
    }
}

# New row order expressed as 0-based indices into $records
# (i.e. $records[0] is the original row 2's data).
# old row 2 -> new row 5, old row 3 -> new row 4, old row 4 -> new row 3,
# old row 5 -> new row 7, old row 6 -> new row 6 (unchanged), old row 7 -> new row 2.
$newOrder = @(5, 2, 1, 0, 4, 3)

# The engine's Hyperlinks.Delete() on a single cell clears every hyperlink
# on the sheet, so drop them all once up front and rebuild in the new order.
$ws.Cells.Item(2, 5).Hyperlinks.Delete()

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $r = 2 + $i
    $rec = $records[$newOrder[$i]]

    $ws.Cells.Item($r, 1).Value = $rec.Title
    $ws.Cells.Item($r, 2).Value = $rec.Timestamp
    $ws.Cells.Item($r, 3).Value = $rec.Distance
    $ws.Cells.Item($r, 4).Value = $rec.Bucket

    $uriCell = $ws.Cells.Item($r, 5)
    $uriCell.Value = $rec.Uri
    $ws.Hyperlinks.Add($uriCell, $rec.Uri) | Out-Null
    $uriCell.Style = "Hyperlink"
}
